# Daily attendance processing - 2026-01-07 07:15:28
# Swap the order of recorders in the "Recorded By" column (G) for every
# session row where both "System" and the instructor email recorded
# attendance: "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
